$d = $word.ActiveDocument

# Paragraph 47 is the empty ListParagraph (numId=4) immediately following the
# "Use-case diagram..." paragraph (46). The edit removes that empty paragraph
# and stamps a collapsed "_GoBack" bookmark at the very end of paragraph 46's
# text (right after its run, before the paragraph mark) - mirroring the
# bookmark Word itself drops at the last edit point.

# Delete the empty paragraph 47 first, so paragraph 46 flows directly into
# what used to be paragraph 48.
$p47 = $d.Paragraphs.Item(47)
$p47.Range.Delete()

$p46 = $d.Paragraphs.Item(46)
$rng46 = $p46.Range
$rng46.MoveEnd(1, -1)
$endPos = $rng46.End

# Temporarily insert a placeholder character right after the visible text so
# we have a genuine mid-run position to anchor the bookmark on (collapsing a
# range exactly on a paragraph-mark boundary does not resolve reliably).
$insertPos = $d.Range($endPos, $endPos)
$insertPos.InsertAfter("X")

$target = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $target)

# Remove the placeholder character again, leaving the bookmark collapsed
# right where the placeholder used to be - i.e. at the true end of paragraph
# 46's text.
$xRange = $d.Range($endPos, $endPos + 1)
$xRange.Delete()
